$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Electrical section: add a new "EEUFC1V331 (Cout)" line item as row 11 ---
# Insert a blank row at 11 (this pushes the old blank "total" row, formerly row 12,
# further down along with everything below it).
$ws.Rows.Item(11).Insert()

# Populate F11 before A11 so that the new shared-string table gets the URL
# registered ahead of the part name (matches upstream ordering).
$ws.Range("F11").Value = "http://goo.gl/rrNymK"
$ws.Range("A11").Value = "EEUFC1V331 (Cout)"
$ws.Range("B11").Value = 0.53
$ws.Range("C11").Value = 1
$ws.Range("D11").Formula = "=SUM(B11*C11,0)"

# The stale blank "total" row that used to be row 12 is now row 12 again (shifted
# down by the insert above); remove it since the new item row replaces its spot.
$ws.Rows.Item(12).Delete()

# --- New "Board Options" sub-section under Mechanical, with its own totals ---
# Row 18 (last Mechanical total row) stays put; insert a new header row at 19,
# then two more blank rows so four blank total rows follow the header (20-23).
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

$ws.Range("A19").Value = "Board Options"
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").HorizontalAlignment = -4108

# The insert at 19 left a stray empty D19 cell behind (content shifted to D20);
# make sure it stays empty.
$ws.Range("D19").ClearContents()

$ws.Range("D20").Formula = "=SUM(B20*C20,0)"
$ws.Range("D21").Formula = "=SUM(B21*C21,0)"
